$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '29.082.03'
$ws.Range("E2").Value = '  +0.16%  '

# Row 3
$ws.Range("D3").Value = '1.835.68'
$ws.Range("E3").Value = '  +0.38%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.16%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '244.77'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.69%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6339'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.22%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.001'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.06%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07554'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.72%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2949'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.31%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '22.97'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.29%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07746'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.68%  '

# Row 12
$ws.Range("D12").Value = '1.835.74'
$ws.Range("E12").Value = '  +0.46%  '

# Row 13
$ws.Range("E13").Value = '  +1.23%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6719'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.51%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '83.35'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.69%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.000009648'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +5.64%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '6.089'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +2.01%  '

# Row 18
$ws.Range("D18").Value = '29.106.48'
$ws.Range("E18").Value = '  +0.26%  '

# Row 19
$ws.Range("E19").Value = '  +2.27%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '226.59'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.98%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.000'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.02%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.210'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.62%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.001'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.13%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '160.82'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.88%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1401'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +3.66%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.556'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.85%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '17.98'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.98%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.504'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.71%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.130'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.20%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.081'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.28%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.203'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.06%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.05394'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +3.58%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.868'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.16%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7466'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.94%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.145'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.34%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.661'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.68%  '

# Row 37
$ws.Range("D37").Value = '1.244.71'
$ws.Range("E37").Value = '  -2.28%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.759'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.38%  '

# Row 39
$ws.Range("E39").Value = '  +0.48%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.635'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +5.06%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9078'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.58%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.001'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.01%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '102.05'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.29%  '

# Row 44
$ws.Range("D44").Value = '1.986.60'
$ws.Range("E44").Value = '  +0.58%  '

# Row 45
$ws.Range("B45").Value = 'BabyDogeCoin'
$ws.Range("C45").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00000000123'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +3.59%  '

# Row 46
$ws.Range("B46").Value = 'Aave'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '64.96'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.65%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5119'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.02%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.4098'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +3.60%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.106'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +3.16%  '

# Row 50
$ws.Range("B50").Value = 'Aptos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '6.785'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.18%  '

# Row 51
$ws.Range("B51").Value = 'RenderToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.653'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.40%  '
